# Apply cryptos list refresh (price / volume updates, plus a Toncoin/Cardano row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.581.63'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.647.86'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.62'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.90'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.628'
$ws.Range('E8').Value = '  +3.09%  '
$ws.Range('E9').Value = '  +3.54%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.81'
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.396'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.69'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '3.125.20'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '65.416.56'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '2.647.37'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.56'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.72'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.43'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.96'
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.64'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.69'
$ws.Range('E26').Value = '  +2.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.58'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '528.12'
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.39'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.36'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '156.99'
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '160.66'
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.61'
$ws.Range('E46').Value = '  -2.41%  '
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0993'
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('D50').Value = '0.0₆0252'
$ws.Range('E50').Value = '  +10.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.74'
$ws.Range('E51').Value = '  -1.11%  '
